# Atualização de bases das ligas, do dia: 03-04-2024 às 22:09
#
# The underlying data source re-sorted several adjacent match rows, which
# results in the entire record (every column except the running index in
# column A) being swapped between each pair of rows below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

$rowPairs = @(
    @(465,466),
    @(485,486),
    @(491,492),
    @(496,497),
    @(509,510),
    @(518,519),
    @(531,532)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    # Snapshot every cell in both rows first so the subsequent writes never
    # read back an already-overwritten value.
    $vals1 = @{}
    $vals2 = @{}
    foreach ($col in $cols) {
        $vals1[$col] = $ws.Range("$col$r1").Value()
        $vals2[$col] = $ws.Range("$col$r2").Value()
    }

    foreach ($col in $cols) {
        $ws.Range("$col$r1").Value = $vals2[$col]
        $ws.Range("$col$r2").Value = $vals1[$col]
    }
}
